# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - row => new F value
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    2  = 1869
    5  = 173
    6  = 2571
    7  = 171
    8  = 89
    10 = 1527
    12 = 44
    13 = 331
    21 = 174
    23 = 1642
    27 = 206
    28 = 298
    29 = 417
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型" (sheet4) - row => new F value
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 1869
    6  = 173
    7  = 2571
    8  = 171
    9  = 89
    11 = 1527
    13 = 44
    14 = 331
    22 = 174
    24 = 1642
    28 = 206
    29 = 298
    30 = 417
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
